$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 367
$ws.Range("I2").Value = 338.85715
$ws.Range("K2").Value = 338.85715
$ws.Range("M2").Value = -225.85715
$ws.Range("H9").Value = 650
$ws.Range("I9").Value = 650
$ws.Range("K9").Value = 650
$ws.Range("M9").Value = -481
$ws.Range("H19").Value = 956.4666999999999
$ws.Range("I19").Value = 609.8
$ws.Range("K19").Value = 609.8
$ws.Range("M19").Value = -434.8
$ws.Range("H29").Value = 2033.2
$ws.Range("I29").Value = 290.5
$ws.Range("K29").Value = 871.5
$ws.Range("M29").Value = -590.5
$ws.Range("H32").Value = 6465.5454
$ws.Range("I32").Value = 7276.5
$ws.Range("J32").Value = 5492.4
$ws.Range("K32").Value = 7276.5
$ws.Range("L32").Value = 5492.4
$ws.Range("M32").Value = -6950.5
$ws.Range("N32").Value = -6144.4
$ws.Range("H43").Value = 6999.5
$ws.Range("I43").Value = 6999
$ws.Range("K43").Value = 6999
$ws.Range("M43").Value = -6930
$ws.Range("H55").Value = 732.875
$ws.Range("J55").Value = 2034.6666
$ws.Range("L55").Value = 2034.6666
$ws.Range("N55").Value = -2462.6666
$ws.Range("H101").Value = 847.8
$ws.Range("J101").Value = 996.6667
$ws.Range("L101").Value = 2990.0001
$ws.Range("N101").Value = -6234.0001
$ws.Range("H116").Value = 15952.857
$ws.Range("I116").Value = 15112.917
$ws.Range("J116").Value = 17072.777
$ws.Range("K116").Value = 15112.917
$ws.Range("L116").Value = 17072.777
$ws.Range("M116").Value = -11670.917
$ws.Range("N116").Value = -23956.777
$ws.Range("H118").Value = 1046.2
$ws.Range("J118").Value = 1966
$ws.Range("L118").Value = 5898
$ws.Range("N118").Value = -9212
$ws.Range("H125").Value = 1720.579
$ws.Range("J125").Value = 1590.1333
$ws.Range("L125").Value = 14311.1997
$ws.Range("N125").Value = -19231.1997
$ws.Range("H127").Value = 1671.1
$ws.Range("I127").Value = 601.5714
$ws.Range("J127").Value = 4166.6665
$ws.Range("K127").Value = 1804.7142
$ws.Range("L127").Value = 12499.9995
$ws.Range("M127").Value = 3155.2858
$ws.Range("N127").Value = -22419.9995
$ws.Range("H131").Value = 4174
$ws.Range("J131").Value = 10526
$ws.Range("L131").Value = 31578
$ws.Range("N131").Value = -41658
$ws.Range("H132").Value = 1818.0416
$ws.Range("I132").Value = 1592.7609
$ws.Range("K132").Value = 4778.2827
$ws.Range("M132").Value = -2248.2827
$ws.Range("H137").Value = 19610710
$ws.Range("I137").Value = 50002004
$ws.Range("K137").Value = 150006012
$ws.Range("M137").Value = -150003462
$ws.Range("H138").Value = 3072.1794
$ws.Range("I138").Value = 1766.3
$ws.Range("J138").Value = 4446.7896
$ws.Range("K138").Value = 5298.9
$ws.Range("L138").Value = 13340.3688
$ws.Range("M138").Value = -158.8999999999996
$ws.Range("N138").Value = -23620.3688
$ws.Range("H141").Value = 1316
$ws.Range("I141").Value = 1088.7561
$ws.Range("K141").Value = 3266.2683
$ws.Range("M141").Value = 1913.7317

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21279610
$ws.Range("I32").Value = 21279610
$ws.Range("K32").Value = 21279610
$ws.Range("M32").Value = -21279323
$ws.Range("H45").Value = 3669.6206
$ws.Range("J45").Value = 6232.1665
$ws.Range("L45").Value = 6232.1665
$ws.Range("N45").Value = -6986.1665
$ws.Range("H74").Value = 18523888
$ws.Range("I74").Value = 30305770
$ws.Range("J74").Value = 9500.429
$ws.Range("K74").Value = 30305770
$ws.Range("L74").Value = 9500.429
$ws.Range("M74").Value = -30304896
$ws.Range("N74").Value = -11248.429
$ws.Range("H77").Value = 18523888
$ws.Range("I77").Value = 30305770
$ws.Range("J77").Value = 9500.429
$ws.Range("K77").Value = 151528850
$ws.Range("L77").Value = 47502.145
$ws.Range("M77").Value = -151524482
$ws.Range("N77").Value = -56238.145
$ws.Range("H132").Value = 4282.56
$ws.Range("I132").Value = 4280.273
$ws.Range("K132").Value = 12840.819
$ws.Range("M132").Value = -10310.819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2566
$ws.Range("H107").Value = 2253.6667
$ws.Range("I107").Value = 1511
$ws.Range("J107").Value = 2625
$ws.Range("K107").Value = 1511
$ws.Range("L107").Value = 2625
$ws.Range("M107").Value = 409
$ws.Range("N107").Value = -6465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1208.3636
$ws.Range("I22").Value = 238.16667
$ws.Range("J22").Value = 2372.6
$ws.Range("K22").Value = 238.16667
$ws.Range("L22").Value = 2372.6
$ws.Range("M22").Value = 111.83333
$ws.Range("N22").Value = -3072.6
$ws.Range("H31").Value = 17527.764
$ws.Range("I31").Value = 1631.6938
$ws.Range("K31").Value = 1631.6938
$ws.Range("M31").Value = -1336.6938
$ws.Range("H34").Value = 17527.764
$ws.Range("I34").Value = 1631.6938
$ws.Range("K34").Value = 1631.6938
$ws.Range("M34").Value = -1429.6938
$ws.Range("H99").Value = 3244
$ws.Range("I99").Value = 3152.3845
$ws.Range("K99").Value = 3152.3845
$ws.Range("M99").Value = -1654.3845
$ws.Range("H107").Value = 899.5625
$ws.Range("I107").Value = 624.3333
$ws.Range("K107").Value = 624.3333
$ws.Range("M107").Value = 1295.6667
$ws.Range("H125").Value = 91663.5
$ws.Range("J125").Value = 91663.5
$ws.Range("L125").Value = 91663.5
$ws.Range("N125").Value = -96583.5
$ws.Range("H126").Value = 3244
$ws.Range("I126").Value = 3152.3845
$ws.Range("K126").Value = 9457.1535
$ws.Range("M126").Value = -6987.1535
$ws.Range("H131").Value = 60000
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2616.4722
$ws.Range("I132").Value = 2269.5588
$ws.Range("K132").Value = 6808.676399999999
$ws.Range("M132").Value = -4278.676399999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 326.66666
$ws.Range("I21").Value = 237.2
$ws.Range("J21").Value = 438.5
$ws.Range("K21").Value = 711.5999999999999
$ws.Range("L21").Value = 1315.5
$ws.Range("M21").Value = -538.5999999999999
$ws.Range("N21").Value = -1661.5
$ws.Range("H86").Value = 2598.6667
$ws.Range("I86").Value = 496.5
$ws.Range("J86").Value = 3649.75
$ws.Range("K86").Value = 1489.5
$ws.Range("L86").Value = 10949.25
$ws.Range("M86").Value = -303.5
$ws.Range("N86").Value = -13321.25
$ws.Range("H89").Value = 2598.6667
$ws.Range("I89").Value = 496.5
$ws.Range("J89").Value = 3649.75
$ws.Range("K89").Value = 4468.5
$ws.Range("L89").Value = 32847.75
$ws.Range("M89").Value = 1459.5
$ws.Range("N89").Value = -44703.75
$ws.Range("H103").Value = 2477.75
$ws.Range("I103").Value = 1424.3334
$ws.Range("J103").Value = 3109.8
$ws.Range("K103").Value = 4273.0002
$ws.Range("L103").Value = 9329.400000000001
$ws.Range("M103").Value = -3394.0002
$ws.Range("N103").Value = -11087.4
$ws.Range("H132").Value = 6639.95
$ws.Range("J132").Value = 7377.6665
$ws.Range("L132").Value = 66398.9985
$ws.Range("N132").Value = -71458.9985

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 917.4545000000001
$ws.Range("I2").Value = 54
$ws.Range("K2").Value = 54
$ws.Range("M2").Value = 59
$ws.Range("H3").Value = 202509.8
$ws.Range("I3").Value = 500425
$ws.Range("J3").Value = 3899.6667
$ws.Range("K3").Value = 500425
$ws.Range("L3").Value = 3899.6667
$ws.Range("M3").Value = -500309
$ws.Range("N3").Value = -4131.6667
$ws.Range("H102").Value = 1962620.9
$ws.Range("I102").Value = 2470287.8
$ws.Range("J102").Value = 4476.857
$ws.Range("K102").Value = 2470287.8
$ws.Range("L102").Value = 4476.857
$ws.Range("M102").Value = -2468665.8
$ws.Range("N102").Value = -7720.857
$ws.Range("H113").Value = 1989.3334
$ws.Range("I113").Value = 1984.5
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 1984.5
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = 185.5
$ws.Range("N113").Value = -6339
$ws.Range("H126").Value = 5773.4
$ws.Range("I126").Value = 2947.25
$ws.Range("J126").Value = 6801.091
$ws.Range("K126").Value = 8841.75
$ws.Range("L126").Value = 20403.273
$ws.Range("M126").Value = -6371.75
$ws.Range("N126").Value = -25343.273

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 14468.389
$ws.Range("I100").Value = 10678.8
$ws.Range("K100").Value = 10678.8
$ws.Range("M100").Value = -10137.8
$ws.Range("H106").Value = 29734.818
$ws.Range("J106").Value = 29734.818
$ws.Range("L106").Value = 29734.818
$ws.Range("N106").Value = -32258.818
$ws.Range("H131").Value = 47999.668
$ws.Range("J131").Value = 48499.5
$ws.Range("L131").Value = 48499.5
$ws.Range("N131").Value = -58579.5
$ws.Range("H136").Value = 6264.857
$ws.Range("I136").Value = 2000.7778
$ws.Range("K136").Value = 6002.3334
$ws.Range("M136").Value = -3452.3334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 403.94116
$ws.Range("I113").Value = 249.2
$ws.Range("K113").Value = 747.5999999999999
$ws.Range("M113").Value = 1422.4
$ws.Range("H122").Value = 2725.1
$ws.Range("I122").Value = 2201.9565
$ws.Range("J122").Value = 4444
$ws.Range("K122").Value = 6605.869499999999
$ws.Range("L122").Value = 13332
$ws.Range("M122").Value = -4155.869499999999
$ws.Range("N122").Value = -18232
